$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (preserves strings such as "1.032" or
# "109.00" that Excel would otherwise silently reinterpret/reformat as numbers),
# then restore the default (unstyled) cell appearance.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- Row 32/33: coin order swapped (ImmutableX now ranks above ARBITRUM) ---
Set-TextValue "B32" "ImmutableX"
Set-TextValue "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "0.7672"
Set-TextValue "E32" "  +2.75%  "
Set-TextValue "B33" "ARBITRUM"
Set-TextValue "C33" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D33" "1.202"
Set-TextValue "E33" "  +5.03%  "

# --- Updated prices / 1h volume percentages for every other row ---
Set-TextValue "D2" "27.522.40"
Set-TextValue "E2" "  +4.09%  "
Set-TextValue "D3" "1.841.76"
Set-TextValue "E3" "  +2.98%  "
Set-TextValue "D4" "1.032"
Set-TextValue "E4" "  +2.91%  "
Set-TextValue "D5" "319.19"
Set-TextValue "E5" "  +4.27%  "
Set-TextValue "D6" "1.027"
Set-TextValue "E6" "  +2.49%  "
Set-TextValue "D7" "0.4371"
Set-TextValue "E7" "  +2.66%  "
Set-TextValue "D8" "0.3737"
Set-TextValue "E8" "  +3.09%  "
Set-TextValue "D9" "0.07377"
Set-TextValue "E9" "  +3.26%  "
Set-TextValue "D10" "0.8742"
Set-TextValue "E10" "  +2.49%  "
Set-TextValue "E11" "  +3.93%  "
Set-TextValue "D12" "1.851.85"
Set-TextValue "E12" "  +2.77%  "
Set-TextValue "D13" "5.490"
Set-TextValue "E13" "  +4.27%  "
Set-TextValue "D14" "6.683"
Set-TextValue "E14" "  +2.88%  "
Set-TextValue "D15" "0.07155"
Set-TextValue "E15" "  +3.47%  "
Set-TextValue "D16" "82.68"
Set-TextValue "E16" "  +3.71%  "
Set-TextValue "D17" "1.031"
Set-TextValue "E17" "  +2.43%  "
Set-TextValue "D18" "0.000008983"
Set-TextValue "E18" "  +2.27%  "
Set-TextValue "E19" "  +2.39%  "
Set-TextValue "D20" "15.39"
Set-TextValue "E20" "  +2.54%  "
Set-TextValue "D21" "27.527.22"
Set-TextValue "E21" "  +4.03%  "
Set-TextValue "D22" "5.257"
Set-TextValue "E22" "  +2.34%  "
Set-TextValue "D23" "11.18"
Set-TextValue "E23" "  +1.15%  "
Set-TextValue "D24" "2.062.43"
Set-TextValue "E24" "  +1.65%  "
Set-TextValue "D25" "157.73"
Set-TextValue "E25" "  +3.66%  "
Set-TextValue "D26" "1.926"
Set-TextValue "E26" "  +5.98%  "
Set-TextValue "E27" "  +2.97%  "
Set-TextValue "D28" "5.254"
Set-TextValue "E28" "  +2.28%  "
Set-TextValue "D29" "1.934"
Set-TextValue "E29" "  +1.14%  "
Set-TextValue "D30" "115.88"
Set-TextValue "E30" "  +1.04%  "
Set-TextValue "E31" "  +2.28%  "
Set-TextValue "D34" "4.494"
Set-TextValue "E34" "  +3.16%  "
Set-TextValue "D35" "2.873"
Set-TextValue "E35" "  +4.56%  "
Set-TextValue "D36" "1.028"
Set-TextValue "E36" "  +2.74%  "
Set-TextValue "D37" "1.143"
Set-TextValue "E37" "  +2.64%  "
Set-TextValue "D38" "0.01973"
Set-TextValue "E38" "  +3.92%  "
Set-TextValue "D39" "0.05249"
Set-TextValue "E39" "  +1.72%  "
Set-TextValue "D40" "0.5168"
Set-TextValue "E40" "  +3.81%  "
Set-TextValue "D41" "2.778"
Set-TextValue "E41" "  +6.57%  "
Set-TextValue "D42" "0.1671"
Set-TextValue "E42" "  +3.03%  "
Set-TextValue "D43" "6.666"
Set-TextValue "E43" "  +4.17%  "
Set-TextValue "D44" "8.529"
Set-TextValue "E44" "  +3.67%  "
Set-TextValue "D45" "109.00"
Set-TextValue "E45" "  +3.29%  "
Set-TextValue "E46" "  +2.86%  "
Set-TextValue "D47" "1.710"
Set-TextValue "E47" "  +4.04%  "
Set-TextValue "D48" "0.4646"
Set-TextValue "E48" "  +3.24%  "
Set-TextValue "D49" "0.06355"
Set-TextValue "D50" "1.888"
Set-TextValue "E50" "  +6.36%  "
Set-TextValue "D51" "39.44"
Set-TextValue "E51" "  +6.70%  "
